$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename vendor values in column C:
#   "HomeSelects"          -> "Home Selects"         (rows 52-59)
#   "HomeSelects Cabinets" -> "Zaca"                  (rows 60-74)
#   "Post Protector Here"  -> "Post Protector-Here"   (row 117)
$lastRow = $ws.UsedRange.Rows.Count

For ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Text
    if ($val -eq "HomeSelects") {
        $cell.Value = "Home Selects"
    } elseif ($val -eq "HomeSelects Cabinets") {
        $cell.Value = "Zaca"
    } elseif ($val -eq "Post Protector Here") {
        $cell.Value = "Post Protector-Here"
    }
}

# Update the view state: scroll so row 106 is the top-left visible row,
# and select G120 as the active cell.
$ws.Range("G120").Select()
$excel.ActiveWindow.ScrollRow = 106
